$d = $word.ActiveDocument

# Remove the "Ver no Jupiter ..." paragraph, the "(c) 2020 ..." paragraph that
# follows it, and the blank paragraph that precedes it (sits right after the
# "LOM3016: ..." requirements line). Walk backwards so deleting a paragraph
# doesn't disturb the indices of paragraphs we still need to examine.

$jupiterText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightText = [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$paras = $d.Paragraphs
$i = $paras.Count
while ($i -ge 1) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($t -eq $jupiterText) {
        $p.Range.Delete()

        $prev = $paras.Item($i - 1)
        if ($prev.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
            $prev.Range.Delete()
        }
    }
    elseif ($t -eq $copyrightText) {
        $p.Range.Delete()
    }

    $i = $i - 1
}

Write-Output "done"
